# Apply updated crypto price/volume values (symbol list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text (price / % change) which must stay
# plain text, matching the workbook's inline-string convention -- force
# Text number format before assigning so Excel does not auto-convert the
# strings into real numbers / percentages.
$targetCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9",
    "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D16", "E16", "D17",
    "E17", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "E24", "D25",
    "E25", "D26", "E26", "D38", "E38", "D39", "E39", "E40", "D41", "E41", "D42", "E42", "D43",
    "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49",
    "D50", "E50", "D51", "E51"
)

foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "304.99"
$ws.Range("E2").Value = "2.57%"
$ws.Range("D3").Value = "44.33"
$ws.Range("E3").Value = "7.36%"
$ws.Range("D4").Value = "5.119"
$ws.Range("E4").Value = "2.25%"
$ws.Range("D5").Value = "0.07807"
$ws.Range("D6").Value = "1.619"
$ws.Range("E6").Value = "3.26%"
$ws.Range("D7").Value = "1.050"
$ws.Range("E7").Value = "12.97%"
$ws.Range("D8").Value = "0.1295"
$ws.Range("E8").Value = "6.58%"
$ws.Range("D9").Value = "0.1864"
$ws.Range("E9").Value = "1.31%"
$ws.Range("D10").Value = "0.09202"
$ws.Range("E10").Value = "3.81%"
$ws.Range("D11").Value = "0.04142"
$ws.Range("E11").Value = "1.79%"
$ws.Range("D12").Value = "0.1045"
$ws.Range("E12").Value = "-0.82%"
$ws.Range("D13").Value = "0.001280"
$ws.Range("E13").Value = "-0.62%"
$ws.Range("D14").Value = "0.005777"
$ws.Range("E14").Value = "-2.27%"
$ws.Range("D16").Value = "3.354"
$ws.Range("E16").Value = "0.30%"
$ws.Range("D17").Value = "4.415"
$ws.Range("E17").Value = "1.21%"
$ws.Range("D19").Value = "0.3363"
$ws.Range("E19").Value = "2.33%"
$ws.Range("D20").Value = "8.052"
$ws.Range("E20").Value = "0.93%"
$ws.Range("D21").Value = "0.1363"
$ws.Range("E21").Value = "-3.87%"
$ws.Range("D22").Value = "0.2809"
$ws.Range("E22").Value = "-5.24%"
$ws.Range("D23").Value = "0.04178"
$ws.Range("E23").Value = "3.27%"
$ws.Range("E24").Value = "0.93%"
$ws.Range("D25").Value = "0.004439"
$ws.Range("E25").Value = "13.74%"
$ws.Range("D26").Value = "0.0001342"
$ws.Range("E26").Value = "9.25%"
$ws.Range("D38").Value = "0.02527"
$ws.Range("E38").Value = "4.43%"
$ws.Range("D39").Value = "0.05339"
$ws.Range("E39").Value = "2.44%"
$ws.Range("E40").Value = "-4.72%"
$ws.Range("D41").Value = "0.007700"
$ws.Range("E41").Value = "-1.20%"
$ws.Range("D42").Value = "0.1369"
$ws.Range("E42").Value = "3.11%"
$ws.Range("D43").Value = "0.007336"
$ws.Range("E43").Value = "-0.36%"
$ws.Range("D44").Value = "0.008323"
$ws.Range("E44").Value = "6.20%"
$ws.Range("D45").Value = "0.3017"
$ws.Range("E45").Value = "1.66%"
$ws.Range("D46").Value = "0.00006677"
$ws.Range("E46").Value = "6.02%"
$ws.Range("D47").Value = "0.00000000746"
$ws.Range("E47").Value = "-0.45%"
$ws.Range("D48").Value = "0.06039"
$ws.Range("E48").Value = "33.71%"
$ws.Range("D49").Value = "0.003978"
$ws.Range("E49").Value = "-5.18%"
$ws.Range("D50").Value = "0.00002088"
$ws.Range("E50").Value = "-0.45%"
$ws.Range("D51").Value = "0.0001989"
$ws.Range("E51").Value = "-0.45%"

# Remove the temporary Text format so the cells end up with no explicit
# style (same as the original workbook), just plain text content.
foreach ($addr in $targetCells) {
    $ws.Range($addr).ClearFormats()
}
